$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: fill in simulations-completed / comments, analyses start date,
# analyses job id, and analyses-completed columns
$ws.Range("F17").Value = "Yes"
$ws.Range("G17").Value = "NA"

# Analyses_start_date (date-only value, formatted mm-dd-yy)
$ws.Range("H17").NumberFormat = "mm-dd-yy"
$ws.Range("H17").Value = [DateTime]::FromOADate(45736)

$ws.Range("I17").Value = "262033 (ac3)"

# Row 16: mark the analyses-completed column as done
$ws.Range("J16").Value = "yes"
$ws.Range("J17").Value = "yes"

# Leave the selection on the newly-entered cell
$ws.Range("J17").Select()
